$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update Price (column D) values that changed but whose row/coin stayed the same ---
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '246.78'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '26.45'
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '5.079'
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.05613'
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '6.501'
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.8130'
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.8441'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '3.749'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004600'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002660'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.008902'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005292'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.002526'

# --- Rows whose Coin/Link/Price/Volume data shifted to a different row (reordering) ---
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.1347'
$ws.Range("E9").Value = '8WazirXWRX'
$ws.Range("B10").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C10").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.03205'
$ws.Range("E10").Value = '9LiechtensteinCryptoassetsExchangeLCX'
$ws.Range("B11").Value = 'BitrueCoin'
$ws.Range("C11").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.02880'
$ws.Range("E11").Value = '10BitrueCoinBTR'
$ws.Range("B12").Value = 'BitMartToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.09413'
$ws.Range("E12").Value = '11BitMartTokenBMX'
$ws.Range("B13").Value = 'BitForexToken'
$ws.Range("C13").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.001513'
$ws.Range("E13").Value = '12BitForexTokenBF'
$ws.Range("B14").Value = 'One'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0005959'
$ws.Range("E14").Value = '13OneONE'
$ws.Range("B15").Value = 'TigerCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.006231'
$ws.Range("E15").Value = '14TigerCashTCH'
$ws.Range("B16").Value = 'LEO'
$ws.Range("C16").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.588'
$ws.Range("E16").Value = '15LEOLEO'
$ws.Range("B17").Value = 'GateToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.016'
$ws.Range("E17").Value = '16GateTokenGT'
$ws.Range("B18").Value = 'BTSEToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.118'
$ws.Range("E18").Value = '17BTSETokenBTSE'
$ws.Range("B19").Value = 'BitpandaEcosystemToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.3157'
$ws.Range("E19").Value = '18BitpandaEcosystemTokenBEST'
$ws.Range("B20").Value = 'MandalaExchangeToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.06976'
$ws.Range("E20").Value = '19MandalaExchangeTokenMDX'
$ws.Range("B41").Value = 'BKEXToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.1348'
$ws.Range("E41").Value = '40BKEXTokenBKKBestin24h'
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.006200'
$ws.Range("E42").Value = '41KickTokenKICK'
